$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string text re-flowed onto multiple indented lines (pretty-printed JSON-ish content).
$newText = @'
questions = [
    {
        "title": "You are defining a data processing pipeline where you want to extract information from several JSON files.Which abstract data type should you use to store the retrieved information?",
        "ques_type": 2,
        "options": [
            "Map",
            "Queue",
            "Tree",
            "Stack"
        ],
        "score": "Map"
    },
    {
        "title": "You are working with time series data where the difference between two data points is one second. After observing the dataset, you notice several missing values in it.Which missing value imputation technique should you use?",
        "ques_type": 2,
        "options": [
            "Imputation with mean",
            "Imputation with mode",
            "Imputation with median",
            "Imputation with previous or next timestamp value"
        ],
        "score": "Imputation with previous or next timestamp value"
    },
    {
        "title": "You are documenting an existing data preprocessing code, and you want to classify all applied changes into two categories: intrarecord structuring and interrecord structuring. You have information about the shape (number of rows x number of columns) of the dataset before and after each modification.Which of the following modifications should you classify as interrecord structuring?",
        "ques_type": 15,
        "options": [
            "(1, 1) -&gt (2, 1)",
            "(1, 1) -&gt (1, 2)",
            "(3, 3) -&gt (6, 6)",
            "(6, 6) -&gt (3, 3)",
            "(10, 8) -&gt (9, 8)"
        ],
        "score": [
            "(3, 3) -&gt (6, 6)",
            "(6, 6) -&gt (3, 3)"
        ]
    },
    {
        "title": "You are preprocessing employee data for your company, and you want to apply syntactic profiling to it. You are currently checking individual employee IDs. The IDs consist of numbers only, with a length of two numbers.What will be the length of the set of valid syntactic values?",
        "ques_type": 2,
        "options": [
            "9",
            "81",
            "99",
            "100"
        ],
        "score": "100"
    }
]
'@

# Row 2 ("A2") held the old shared-string text; it goes away entirely, collapsing the
# used range back down to just A1.
$ws.Range("A2").ClearContents()

# Row 1 ("A1") used to hold a bare 0 formatted with a bold font + thin box border +
# centered alignment (cellXf index 1). That formatting is dropped -> back to the
# default "Normal" style, and the cell now carries the (reformatted) text that used
# to live in A2.
$r1 = $ws.Range("A1")
$r1.Style = "Normal"
$r1.Value = $newText

# Typing a value with embedded line breaks makes Excel auto-grow the row height;
# AutoFit puts row 1 back to the sheet's default (no explicit/custom height).
$ws.Rows(1).AutoFit()
